$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3853.0667
$ws.Range("I40").Value = 2666.75
$ws.Range("J40").Value = 5208.857
$ws.Range("K40").Value = 2666.75
$ws.Range("L40").Value = 5208.857
$ws.Range("M40").Value = -2491.75
$ws.Range("N40").Value = -5558.857

$ws.Range("H51").Value = 4999
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 4999
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 4999
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -5967

$ws.Range("H86").Value = 3597.4
$ws.Range("I86").Value = 3769.8572
$ws.Range("J86").Value = 3195
$ws.Range("K86").Value = 3769.8572
$ws.Range("L86").Value = 3195
$ws.Range("M86").Value = -2646.8572
$ws.Range("N86").Value = -5441

$ws.Range("H89").Value = 3597.4
$ws.Range("I89").Value = 3769.8572
$ws.Range("J89").Value = 3195
$ws.Range("K89").Value = 18849.286
$ws.Range("L89").Value = 15975
$ws.Range("M89").Value = -13233.286
$ws.Range("N89").Value = -27207

$ws.Range("H141").Value = 2481.342
$ws.Range("I141").Value = 1872.7297
$ws.Range("K141").Value = 5618.189100000001
$ws.Range("M141").Value = -438.1891000000005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 350
$ws.Range("J4").Value = 400
$ws.Range("L4").Value = 400
$ws.Range("N4").Value = -632

$ws.Range("H32").Value = 11499164
$ws.Range("I32").Value = 12350190
$ws.Range("K32").Value = 12350190
$ws.Range("M32").Value = -12349903

$ws.Range("H74").Value = 1303.7333
$ws.Range("I74").Value = 836.2
$ws.Range("K74").Value = 836.2
$ws.Range("M74").Value = 37.79999999999995

$ws.Range("H77").Value = 1303.7333
$ws.Range("I77").Value = 836.2
$ws.Range("K77").Value = 4181
$ws.Range("M77").Value = 187

$ws.Range("H133").Value = 69964.516
$ws.Range("J133").Value = 69964.516
$ws.Range("L133").Value = 69964.516
$ws.Range("N133").Value = -75024.516

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 80065
$ws.Range("J140").Value = 80065
$ws.Range("L140").Value = 80065
$ws.Range("N140").Value = -90425

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 240.78572
$ws.Range("I7").Value = 112.14286
$ws.Range("J7").Value = 369.42856
$ws.Range("K7").Value = 112.14286
$ws.Range("L7").Value = 369.42856
$ws.Range("M7").Value = 0.8571400000000011
$ws.Range("N7").Value = -595.4285600000001

$ws.Range("H31").Value = 3014.3125
$ws.Range("I31").Value = 2827.6924
$ws.Range("J31").Value = 3823
$ws.Range("K31").Value = 2827.6924
$ws.Range("L31").Value = 3823
$ws.Range("M31").Value = -2532.6924
$ws.Range("N31").Value = -4413

$ws.Range("H34").Value = 3014.3125
$ws.Range("I34").Value = 2827.6924
$ws.Range("J34").Value = 3823
$ws.Range("K34").Value = 2827.6924
$ws.Range("L34").Value = 3823
$ws.Range("M34").Value = -2625.6924
$ws.Range("N34").Value = -4227

$ws.Range("H132").Value = 1474.1666
$ws.Range("I132").Value = 1443.7778
$ws.Range("J132").Value = 1565.3334
$ws.Range("K132").Value = 4331.3334
$ws.Range("L132").Value = 4696.0002
$ws.Range("M132").Value = -1801.3334
$ws.Range("N132").Value = -9756.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 29735866
$ws.Range("I4").Value = 37444916
$ws.Range("K4").Value = 112334748
$ws.Range("M4").Value = -112334636

$ws.Range("H131").Value = 2445.5652
$ws.Range("I131").Value = 1487.45
$ws.Range("J131").Value = 8833
$ws.Range("K131").Value = 4462.35
$ws.Range("L131").Value = 26499
$ws.Range("M131").Value = 577.6499999999996
$ws.Range("N131").Value = -36579

$ws.Range("H134").Value = 9037.375
$ws.Range("I134").Value = 2659.8
$ws.Range("K134").Value = 7979.400000000001
$ws.Range("M134").Value = -2909.400000000001

$ws.Range("H137").Value = 2017.3636
$ws.Range("J137").Value = 2266.6
$ws.Range("L137").Value = 6799.799999999999
$ws.Range("N137").Value = -16999.8

$ws.Range("H139").Value = 4153.579
$ws.Range("I139").Value = 2422.7144
$ws.Range("K139").Value = 7268.1432
$ws.Range("M139").Value = -2128.1432

$ws.Range("H140").Value = 4504.6665
$ws.Range("I140").Value = 4504.6665
$ws.Range("K140").Value = 13513.9995
$ws.Range("M140").Value = -8333.999500000002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws.Range("H46").Value = 10000
$ws.Range("I46").Value = 10000
$ws.Range("K46").Value = 10000
$ws.Range("M46").Value = -9844

$ws.Range("H57").Value = 4999
$ws.Range("I57").Value = 4999
$ws.Range("K57").Value = 4999
$ws.Range("M57").Value = -4179

$ws.Range("H80").Value = 94927
$ws.Range("I80").Value = 223634
$ws.Range("J80").Value = 2993.4285
$ws.Range("K80").Value = 223634
$ws.Range("L80").Value = 2993.4285
$ws.Range("M80").Value = -222636
$ws.Range("N80").Value = -4989.4285

$ws.Range("H83").Value = 94927
$ws.Range("I83").Value = 223634
$ws.Range("J83").Value = 2993.4285
$ws.Range("K83").Value = 1118170
$ws.Range("L83").Value = 14967.1425
$ws.Range("M83").Value = -1113178
$ws.Range("N83").Value = -24951.1425

$ws.Range("H107").Value = 1087.2307
$ws.Range("I107").Value = 1283
$ws.Range("J107").Value = 646.75
$ws.Range("K107").Value = 1283
$ws.Range("L107").Value = 646.75
$ws.Range("M107").Value = 637
$ws.Range("N107").Value = -4486.75

$ws.Range("H126").Value = 4007.1667
$ws.Range("I126").Value = 4007.1667
$ws.Range("K126").Value = 12021.5001
$ws.Range("M126").Value = -9551.500100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 56166
$ws.Range("J36").Value = 56166
$ws.Range("L36").Value = 56166
$ws.Range("N36").Value = -57290

$ws.Range("H131").Value = 69799
$ws.Range("J131").Value = 69799
$ws.Range("L131").Value = 69799
$ws.Range("N131").Value = -79879

$ws.Range("H132").Value = 7158.0264
$ws.Range("I132").Value = 4387.375
$ws.Range("K132").Value = 13162.125
$ws.Range("M132").Value = -10632.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 36690.812
$ws.Range("J96").Value = 7145.25
$ws.Range("L96").Value = 7145.25
$ws.Range("N96").Value = -9891.25

$ws.Range("H136").Value = 1359.5
$ws.Range("I136").Value = 943.65216
$ws.Range("J136").Value = 1922.1177
$ws.Range("K136").Value = 2830.95648
$ws.Range("L136").Value = 5766.3531
$ws.Range("M136").Value = -280.9564799999998
$ws.Range("N136").Value = -10866.3531
